$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.552.88'
$ws.Range("E2").Value = '  +2.53%  '
$ws.Range("D3").Value = '1.642.55'
$ws.Range("E3").Value = '  +4.35%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.33'
$ws.Range("E5").Value = '  +2.95%  '
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3783'
$ws.Range("E7").Value = '  +0.95%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '53.13'
$ws.Range("E8").Value = '  +6.47%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3690'
$ws.Range("E9").Value = '  +3.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.285'
$ws.Range("E10").Value = '  +5.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08205'
$ws.Range("E11").Value = '  +3.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9993'
$ws.Range("E12").Value = '  -0.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.27'
$ws.Range("E13").Value = '  +6.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.681'
$ws.Range("E14").Value = '  +4.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001290'
$ws.Range("E15").Value = '  +5.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.478'
$ws.Range("E16").Value = '  +2.90%  '
$ws.Range("D17").Value = '1.641.61'
$ws.Range("E17").Value = '  +4.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.10'
$ws.Range("E18").Value = '  +3.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06962'
$ws.Range("E19").Value = '  +3.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.48'
$ws.Range("E20").Value = '  +4.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.604'
$ws.Range("E21").Value = '  +4.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9985'
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("D23").Value = '23.560.94'
$ws.Range("E23").Value = '  +2.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.02'
$ws.Range("E24").Value = '  +3.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.138'
$ws.Range("E25").Value = '  +12.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.418'
$ws.Range("E26").Value = '  +2.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.48'
$ws.Range("E27").Value = '  +4.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.90'
$ws.Range("E28").Value = '  +3.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.345'
$ws.Range("E29").Value = '  +3.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '136.53'
$ws.Range("E30").Value = '  +4.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.427'
$ws.Range("E31").Value = '  +4.20%  '
$ws.Range("E32").Value = '  +5.39%  '
$ws.Range("D33").Value = '1.815.53'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9784'
$ws.Range("E34").Value = '  +5.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02833'
$ws.Range("E35").Value = '  +7.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.45'
$ws.Range("E36").Value = '  +5.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.07496'
$ws.Range("E37").Value = '  +2.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.245'
$ws.Range("E38").Value = '  +4.90%  '
$ws.Range("E39").Value = '  +3.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08868'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.403'
$ws.Range("E41").Value = '  +4.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7194'
$ws.Range("E42").Value = '  +5.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.69'
$ws.Range("E43").Value = '  +7.66%  '
$ws.Range("E44").Value = '  +10.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6658'
$ws.Range("E45").Value = '  +5.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.371'
$ws.Range("E46").Value = '  +5.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.045'
$ws.Range("E47").Value = '  +1.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9977'
$ws.Range("E48").Value = '  -0.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08079'
$ws.Range("E49").Value = '  +3.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.54'
$ws.Range("E50").Value = '  +1.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.219'
$ws.Range("E51").Value = '  +3.31%  '
